$d = $word.ActiveDocument

# --- wdFind* constants (not pre-defined in this host, declare locally) ---
$wdReplaceNone = 0
$wdReplaceOne  = 1
$wdReplaceAll  = 2
$wdFindStop    = 0

# ------------------------------------------------------------------
# 1) Reference 5 (WHO Classification of Tumours) — three runs carrying
#    " ... vol. 11). " + "https://publications.iarc.fr" + ".  " get
#    collapsed into a single run with the concatenated text.
# ------------------------------------------------------------------
$ref5Old = " WHO Classification of Tumours Editorial Board. Haematolymphoid tumours. Lyon (France): International Agency for Research on Cancer; forthcoming. (WHO classification of tumours series, 5th ed.; vol. 11). "
$ref5Old += "https://publications.iarc.fr"
$ref5Old += ".  "

$ref5New = $ref5Old

$r5 = $d.Content
$r5.Find.Execute($ref5Old, $false, $false, $false, $false, $false, $true, $wdFindStop, $false, $ref5New, $wdReplaceOne)

# ------------------------------------------------------------------
# 2) Reference 12 — Schuurhuis GJ 2018 Blood 131(12):1275-91  ->
#    Heuser M 2021 Blood 138(26):2753-67. Each piece below is the
#    full text of its own (isolated) run, so replacing it in place
#    keeps the surrounding run/formatting boundaries untouched.
# ------------------------------------------------------------------

# Anchor on the unique start of the citation so later generic
# fragments ("2018;", "131", etc. which also occur elsewhere in the
# references list) only get touched inside this citation's vicinity.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Schuurhuis GJ, et al. Minimal/measurable residual disease")
$scopeStart = $anchor.Start - 1
$scope = $d.Range($scopeStart, $scopeStart + 400)

$scope.Find.Execute(
    " Schuurhuis GJ, et al. Minimal/measurable residual disease in AML: a consensus document from the European LeukemiaNet MRD Working Party. ",
    $false, $false, $false, $false, $false, $true, $wdFindStop, $false,
    " Heuser M, et al. 2021 Update on MRD in acute myeloid leukemia: a consensus document from the European LeukemiaNet MRD Working Party. ",
    $wdReplaceOne)

$scope = $d.Range($scopeStart, $scopeStart + 400)
$scope.Find.Execute(" 2018; ", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, " 2021; ", $wdReplaceOne)

$scope = $d.Range($scopeStart, $scopeStart + 400)
$scope.Find.Execute("131", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, "138", $wdReplaceOne)

$scope = $d.Range($scopeStart, $scopeStart + 400)
$scope.Find.Execute("(12): 1275-91.  ", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, "(26): 2753-67.  ", $wdReplaceOne)

# ------------------------------------------------------------------
# 3) Two additional empty (small, hidden-style) paragraphs right
#    after the closing references table, before the trailing
#    paragraph that already carries that formatting.
# ------------------------------------------------------------------
$tbl = $d.Tables.Item($d.Tables.Count)
$afterTable = $tbl.Range.End
$insertPoint = $d.Range($afterTable, $afterTable)

for ($i = 0; $i -lt 2; $i++) {
    $para = $d.Paragraphs.Add($insertPoint)
    $para.Range.Font.Name = "Calibri"
    $para.Range.Font.NameAscii = "Calibri"
    $para.Range.Font.Size = 3
    $insertPoint = $d.Range($para.Range.End, $para.Range.End)
}
